$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.343.67"
$ws.Range("E2").Value = "  +4.49%  "
$ws.Range("D3").Value = "'2.367.98"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'310.14"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'108.32"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "'41.13"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "'8.48"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D15").Value = "'2.729.65"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "'15.22"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").Value = "'2.357.95"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "'45.306.26"
$ws.Range("E18").Value = "  +4.63%  "
$ws.Range("D19").Value = "'14.44"
$ws.Range("E19").Value = "  +11.24%  "
$ws.Range("D20").Value = "'7.33"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "'3.49"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "'260.60"
$ws.Range("E24").Value = "  -2.99%  "
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'7.35"
$ws.Range("E28").Value = "  -3.20%  "
$ws.Range("D29").Value = "'2.35"
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("D30").Value = "'0.0971"
$ws.Range("E30").Value = "  +9.80%  "
$ws.Range("D31").Value = "'22.36"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("D33").Value = "'169.19"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "'2.95"
$ws.Range("E34").Value = "  +6.20%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("E36").Value = "  +4.50%  "
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("E38").Value = "  +5.05%  "
$ws.Range("D39").Value = "'3.93"
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("D42").Value = "'99.55"
$ws.Range("E42").Value = "  -5.40%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").Value = "'69.59"
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("D45").Value = "'12.97"
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "'81.02"
$ws.Range("E47").Value = "  +5.34%  "
$ws.Range("D48").Value = "'112.15"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").Value = "'5.54"
$ws.Range("E49").Value = "  +4.13%  "
$ws.Range("D50").Value = "'9.22"
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("D51").Value = "'1.668.13"
$ws.Range("E51").Value = "  +0.38%  "
